$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3439717292785645
$ws.Range("B1").Value = 0.4329738914966583
$ws.Range("C1").Value = 0.6292267441749573
$ws.Range("D1").Value = 2.490821123123169
$ws.Range("E1").Value = 5.614772319793701
